# Update F-column ("想去人数" / interest counts) across all four sheets
# to match the regenerated gh-pages data snapshot (commit 456a3b4).
$wb = $excel.ActiveWorkbook

# Exhibition (展览)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 580
$ws.Range("F10").Value = 940
$ws.Range("F14").Value = 92
$ws.Range("F16").Value = 876
$ws.Range("F17").Value = 1779
$ws.Range("F18").Value = 3542
$ws.Range("F19").Value = 1043
$ws.Range("F21").Value = 2454
$ws.Range("F22").Value = 676
$ws.Range("F23").Value = 38
$ws.Range("F24").Value = 3387
$ws.Range("F26").Value = 821
$ws.Range("F28").Value = 2049
$ws.Range("F33").Value = 141
$ws.Range("F35").Value = 1221
$ws.Range("F36").Value = 1878
$ws.Range("F37").Value = 456
$ws.Range("F40").Value = 241
$ws.Range("F43").Value = 67

# Performance (演出)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 6

# Local life (本地生活)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 212

# All types (全部类型)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 580
$ws.Range("F8").Value = 940
$ws.Range("F14").Value = 92
$ws.Range("F15").Value = 876
$ws.Range("F16").Value = 1779
$ws.Range("F17").Value = 3542
$ws.Range("F18").Value = 1043
$ws.Range("F21").Value = 2454
$ws.Range("F23").Value = 38
$ws.Range("F24").Value = 3387
$ws.Range("F26").Value = 821
$ws.Range("F27").Value = 6
$ws.Range("F29").Value = 2049
$ws.Range("F38").Value = 141
$ws.Range("F40").Value = 1221
$ws.Range("F41").Value = 1878
$ws.Range("F43").Value = 456
$ws.Range("F45").Value = 241
$ws.Range("F48").Value = 67

